# The first two graded-submission rows (row 2 and row 3 of the data, i.e.
# aid 2 and aid 3) are moved down to the bottom of the list, after the last
# existing row (old row 15 / aid 15). Every other row shifts up by two rows
# to fill the gap. Net effect on the data block A2:J15: a cyclic shift where
# old rows 4..15 become new rows 2..13, and old rows 2..3 become new rows
# 14..15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A2:J15")
$data = $range.Value2

$rowCount = $data.GetLength(0)
$colCount = $data.GetLength(1)

# Build the new row order (1-based, relative to the top of $data):
# rows 3..14 first (old rows 4..15), then rows 1..2 (old rows 2..3).
$order = @()
for ($i = 3; $i -le $rowCount; $i++) { $order += $i }
$order += 1
$order += 2

$reordered = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $order[$i]
    for ($j = 1; $j -le $colCount; $j++) {
        $reordered[$i, $j - 1] = $data[$srcRow, $j]
    }
}

$range.Value2 = $reordered
